# Insert a new weekly price record as row 276 on the "Espinaca" sheet.
# This pushes the existing rows 276:355 down to 277:356 (dimension grows
# from A1:R355 to A1:R356) and populates the newly opened row 276 with a
# fresh data record for 2021-09-24 (serial 44463).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 276, shifting 276:355 -> 277:356.
$ws.Rows("276:276").Insert()

# Populate the new row 276 with the new record's values.
$ws.Range("A276").Value = 6
$ws.Range("B276").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C276").Value = "Metropolitana"
$ws.Range("D276").Value = 44463
$ws.Range("E276").Value = 13
$ws.Range("F276").Value = 100112012
$ws.Range("G276").Value = "Espinaca"
$ws.Range("H276").Value = "Sin especificar"
$ws.Range("I276").Value = "Primera"
$ws.Range("J276").Value = 400
$ws.Range("K276").Value = 4500
$ws.Range("L276").Value = 5000
$ws.Range("M276").Value = 4788
$ws.Range("N276").Value = "$/cuna 10 kilos"
$ws.Range("O276").Value = "Región Metropolitana"
$ws.Range("P276").Value = 479
$ws.Range("Q276").Value = 10
$ws.Range("R276").Value = "Hortaliza"
